$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Entry")

# The TRANSECT column (F) for rows 55-66 currently reads "Sandpiper".
# Rename it to the fuller "Sandpiper Beach" site/transect name, matching
# the cleaned-up kelp data.
$range = $ws.Range("F55:F66")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq "Sandpiper") {
        $cell.Value = "Sandpiper Beach"
    }
}
